$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.183.35"
$ws.Range("E2").Value = "  +1.82%  "

# Row 3
$ws.Range("D3").Value = "3.922.87"
$ws.Range("E3").Value = "  +0.51%  "

# Row 4
$c = $ws.Range("D4")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$c = $ws.Range("D5")
$c.Value = "'484.38"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +3.76%  "

# Row 6
$c = $ws.Range("D6")
$c.Value = "'146.45"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.71%  "

# Row 7
$c = $ws.Range("D7")
$c.Value = "'0.628"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.14%  "

# Row 8
$ws.Range("E8").Value = "  -0.20%  "

# Row 9
$c = $ws.Range("D9")
$c.Value = "'0.729"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.35%  "

# Row 10
$ws.Range("E10").Value = "  +2.05%  "

# Row 11
$c = $ws.Range("D11")
$c.Value = "'0.0000357"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +4.40%  "

# Row 12
$c = $ws.Range("D12")
$c.Value = "'42.63"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.43%  "

# Row 13
$c = $ws.Range("D13")
$c.Value = "'10.61"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.68%  "

# Row 14
$ws.Range("D14").Value = "4.538.75"
$ws.Range("E14").Value = "  +0.26%  "

# Row 15
$c = $ws.Range("D15")
$c.Value = "'14.85"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.11%  "

# Row 16
$ws.Range("D16").Value = "3.951.09"
$ws.Range("E16").Value = "  +0.41%  "

# Row 17
$ws.Range("E17").Value = "  -0.15%  "

# Row 18
$c = $ws.Range("D18")
$c.Value = "'19.88"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.62%  "

# Row 19
$ws.Range("E19").Value = "  -2.14%  "

# Row 20
$ws.Range("D20").Value = "68.320.64"
$ws.Range("E20").Value = "  +1.65%  "

# Row 21
$c = $ws.Range("D21")
$c.Value = "'447.71"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +3.44%  "

# Row 22
$ws.Range("E22").Value = "  -0.39%  "

# Row 23
$c = $ws.Range("D23")
$c.Value = "'3.35"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.07%  "

# Row 24
$c = $ws.Range("D24")
$c.Value = "'89.04"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.48%  "

# Row 25
$c = $ws.Range("D25")
$c.Value = "'11.51"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +14.30%  "

# Row 26
$c = $ws.Range("D26")
$c.Value = "'10.91"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +12.95%  "

# Row 27
$c = $ws.Range("D27")
$c.Value = "'3.61"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.37%  "

# Row 28
$c = $ws.Range("D28")
$c.Value = "'38.90"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.79%  "

# Row 29
$c = $ws.Range("D29")
$c.Value = "'5.85"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +3.41%  "

# Row 30
$c = $ws.Range("D30")
$c.Value = "'13.42"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.65%  "

# Row 31
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D31")
$c.Value = "'0.131"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.52%  "

# Row 32
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D32")
$c.Value = "'688.92"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -6.93%  "

# Row 33
$c = $ws.Range("D33")
$c.Value = "'2.87"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +3.02%  "

# Row 34
$ws.Range("D34").Value = "0.0₃0930"
$ws.Range("E34").Value = "  +18.94%  "

# Row 35
$c = $ws.Range("D35")
$c.Value = "'41.88"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -4.20%  "

# Row 36
$c = $ws.Range("D36")
$c.Value = "'59.13"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +2.01%  "

# Row 37
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D37")
$c.Value = "'0.150"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -5.12%  "

# Row 38
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D38")
$c.Value = "'5.66"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +5.71%  "

# Row 39
$ws.Range("E39").Value = "  +0.00%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D40")
$c.Value = "'0.0478"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.24%  "

# Row 41
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D41")
$c.Value = "'2.87"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +16.07%  "

# Row 42
$c = $ws.Range("D42")
$c.Value = "'0.366"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +9.04%  "

# Row 43
$ws.Range("E43").Value = "  -5.49%  "

# Row 44
$ws.Range("E44").Value = "  +6.81%  "

# Row 45
$ws.Range("E45").Value = "  +0.53%  "

# Row 46
$ws.Range("E46").Value = "  +0.04%  "

# Row 47
$ws.Range("E47").Value = "  -0.27%  "

# Row 48
$c = $ws.Range("D48")
$c.Value = "'2.14"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -2.14%  "

# Row 49
$c = $ws.Range("D49")
$c.Value = "'146.23"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.86%  "

# Row 50
$ws.Range("E50").Value = "  -1.48%  "

# Row 51
$c = $ws.Range("D51")
$c.Value = "'2.84"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.52%  "
